$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("AN2").Value = 6
$ws.Range("AO2").Value = 11
$ws.Range("G2").Value = 3.3
$ws.Range("M2").Value = 1.17
$ws.Range("N2").Value = 5
$ws.Range("W2").Value = 6.5
$ws.Range("X2").Value = 1.11

# Row 3
$ws.Range("AO3").Value = 29
$ws.Range("AP3").Value = 21
$ws.Range("I3").Value = 6
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("U3").Value = 4.2
$ws.Range("V3").Value = 1.23

# Row 5
$ws.Range("S5").Value = 2.2
$ws.Range("T5").Value = 1.65

# Row 6
$ws.Range("AI6").Value = 6.5
$ws.Range("AL6").Value = 67
$ws.Range("AM6").Value = 451
$ws.Range("AQ6").Value = 21
$ws.Range("G6").Value = 3.7
$ws.Range("I6").Value = 2.25
$ws.Range("K6").Value = 1.95
$ws.Range("L6").Value = 3

# Row 7
$ws.Range("AA7").Value = 1.91
$ws.Range("AB7").Value = 1.7
$ws.Range("AC7").Value = 9.25
$ws.Range("AD7").Value = 6.8
$ws.Range("AE7").Value = 9.5
$ws.Range("AF7").Value = 7.4
$ws.Range("AG7").Value = 10.25
$ws.Range("AH7").Value = 27
$ws.Range("AI7").Value = 18.5
$ws.Range("AJ7").Value = 13
$ws.Range("AK7").Value = 25
$ws.Range("AL7").Value = 110
$ws.Range("AM7").Value = 800
$ws.Range("AO7").Value = 80
$ws.Range("AP7").Value = 35
$ws.Range("AQ7").Value = 300
$ws.Range("AR7").Value = 120
$ws.Range("AS7").Value = 100
$ws.Range("G7").Value = 1.21
$ws.Range("H7").Value = 6
$ws.Range("I7").Value = 10.5
$ws.Range("J7").Value = 1.57
$ws.Range("K7").Value = 2.75
$ws.Range("L7").Value = 8.25
$ws.Range("S7").Value = 1.4
$ws.Range("T7").Value = 2.52
$ws.Range("W7").Value = 1.98
$ws.Range("X7").Value = 1.65

# Row 8
$ws.Range("AC8").Value = 9
$ws.Range("AD8").Value = 17
$ws.Range("AE8").Value = 13
$ws.Range("AF8").Value = 41
$ws.Range("AM8").Value = 351
$ws.Range("AO8").Value = 9.5
$ws.Range("AR8").Value = 19
$ws.Range("G8").Value = 3.5
$ws.Range("H8").Value = 3.2
$ws.Range("I8").Value = 2.2
$ws.Range("L8").Value = 3
$ws.Range("N8").Value = 7.5

# Row 12
$ws.Range("AA12").Value = 2.6
$ws.Range("AD12").Value = 4.9
$ws.Range("AE12").Value = 9.5
$ws.Range("AF12").Value = 6.9
$ws.Range("AI12").Value = 7.1
$ws.Range("AJ12").Value = 9.75
$ws.Range("AR12").Value = 200
$ws.Range("G12").Value = 1.28
$ws.Range("H12").Value = 4.6
$ws.Range("I12").Value = 12
$ws.Range("K12").Value = 2.35
$ws.Range("L12").Value = 9.75
$ws.Range("N12").Value = 7.1
$ws.Range("O12").Value = 1.32
$ws.Range("P12").Value = 3.1
$ws.Range("S12").Value = 1.93
$ws.Range("T12").Value = 1.78
$ws.Range("W12").Value = 3.25
$ws.Range("Y12").Value = 1.38
$ws.Range("Z12").Value = 2.8

# Row 14
$ws.Range("AC14").Value = 6
$ws.Range("I14").Value = 3.75
$ws.Range("J14").Value = 2.75
$ws.Range("L14").Value = 4.5
$ws.Range("N14").Value = 8

# Row 16
$ws.Range("AG16").Value = 51
$ws.Range("AH16").Value = 51
$ws.Range("AJ16").Value = 8.5
$ws.Range("AK16").Value = 19
$ws.Range("H16").Value = 4.33
$ws.Range("N16").Value = 12

# Row 20
$ws.Range("AN20").Value = 6.5
$ws.Range("G20").Value = 4.75
$ws.Range("I20").Value = 1.73
$ws.Range("L20").Value = 2.38
$ws.Range("M20").Value = 1.06
$ws.Range("N20").Value = 10
$ws.Range("O20").Value = 1.3
$ws.Range("P20").Value = 3.4
$ws.Range("S20").Value = 2.05
$ws.Range("T20").Value = 1.8
$ws.Range("W20").Value = 3.5
$ws.Range("X20").Value = 1.29

# Row 21
$ws.Range("M21").Value = 1.04
$ws.Range("N21").Value = 13
$ws.Range("S21").Value = 1.73
$ws.Range("T21").Value = 2.08

# Row 22
$ws.Range("AE22").Value = 11
$ws.Range("I22").Value = 2.25
$ws.Range("N22").Value = 8.5
$ws.Range("S22").Value = 2.15
$ws.Range("T22").Value = 1.67
$ws.Range("Y22").Value = 1.5
$ws.Range("Z22").Value = 2.5

# Row 23
$ws.Range("N23").Value = 10
$ws.Range("W23").Value = 3.5
$ws.Range("X23").Value = 1.29

# Row 24
$ws.Range("AC24").Value = 8
$ws.Range("AD24").Value = 10
$ws.Range("AG24").Value = 15
$ws.Range("AH24").Value = 23
$ws.Range("AN24").Value = 12
$ws.Range("AP24").Value = 13
$ws.Range("AR24").Value = 29
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = 3.4
$ws.Range("I24").Value = 3.6
$ws.Range("J24").Value = 2.63
$ws.Range("S24").Value = 1.88
$ws.Range("T24").Value = 1.98

# Row 26
$ws.Range("AC26").Value = 11
$ws.Range("G26").Value = 3.3
$ws.Range("I26").Value = 2.15
$ws.Range("M26").Value = 1.05
$ws.Range("N26").Value = 11
